# BOM WT SU assembly - add brake disc rows (images des assemblages)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row before the "Rear hubs" section header (old row 17) for
#    the Front Brake Disc part, matching the style/format of the row above.
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).Insert()
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Append a new row at the end of the sheet (old last row 27, now row 28)
#    for the Rear Brake Disc part, matching the style/format of the row above.
# ---------------------------------------------------------------------------
$ws.Range("A28:G28").Copy()
$ws.Range("A29:G29").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Fill in the values. "Rear Brake Disc" is entered before "Front Brake
#    Disc" so the shared-string table order matches the source edit.
# ---------------------------------------------------------------------------
$ws.Range("C29").Value2 = "Rear Brake Disc"
$ws.Range("D29").Value2 = "b"
$ws.Range("F29").Value2 = 2
$ws.Range("G29").Value2 = "WT_03010"

$ws.Range("C17").Value2 = "Front Brake Disc"
$ws.Range("D17").Value2 = "b"
$ws.Range("F17").Value2 = 2
$ws.Range("G17").Value2 = "WT_02011"

$excel.CutCopyMode = 0
